$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENTRADAS")

# Existing values are in row 1 (A1=G101, B1=E102, C1=S104).
# New layout: header row, then one ID per row, with date/time on the last (matched) row.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "FECHA"
$ws.Range("C1").Value = "HORA"

$ws.Range("A2").Value = "G101"
$ws.Range("A3").Value = "E102"

$ws.Range("A4").Value = "S104"
# Leading apostrophe forces Excel to store this as literal text instead of
# auto-converting the ISO-looking string into a date serial number.
$ws.Range("B4").Value = "'2024-12-16"
$ws.Range("C4").Value = "11:51:27"
